$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Solitario",    "Annabmota", 209, 9,  "2025-11-06 19:21:47"),
    @("Multijugador",  "Jimena",    283, 10, "2025-11-06 19:34:53"),
    @("Solitario",    "Annita",    834, 5,  "2025-11-06 19:40:53"),
    @("Solitario",    "Prueba",    516, 5,  "2025-11-06 19:53:02")
)

$rowIndex = 4
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
